$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Builder" for the new data block: row 3 is filled with a constant (1),
# row 4 is filled with an incrementing 1..20 sequence, both spanning
# columns A (1) through T (20).
$columnCount = 20

for ($col = 1; $col -le $columnCount; $col++) {
    $ws.Cells.Item(3, $col).Value = 1
    $ws.Cells.Item(4, $col).Value = $col
}

# Scroll the view so column H is the left-most visible column, matching
# the author's on-screen position while reviewing the newly built data.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1

# Move/extend the selection to the freshly built row.
$ws.Range("A4:T4").Select()
